# Update countries & provincias Spain
#
# 1) Fix the Fiyi / Dominica ordering: Dominica should appear on the row that
#    previously showed Fiyi, and Fiyi moves to the row right below it
#    (the two countries have identical totals, so only the labels swap).
# 2) Bump the "Datos actualizados" timestamp in A1.
# 3) Refresh the daily COVID figures for a handful of countries
#    (Kazajistan, Haiti, Montenegro, Mongolia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the Fiyi / Dominica country labels -----------------------------
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# --- 2. Update the "last updated" timestamp ---------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 06:28"

# --- 3. Refresh country statistics -------------------------------------------

# Kazajistan (row 38): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B38").Value = 47171
$ws.Range("C38").Value = 1452
$ws.Range("D38").Value = 27030
$ws.Range("E38").Value = 19953

# Haiti (row 84): Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes
$ws.Range("B84").Value = 6294
$ws.Range("C84").Value = 64
$ws.Range("D84").Value = 1408
$ws.Range("E84").Value = 4773
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 113

# Montenegro (row 147): Recuperados, Casos criticos, Muertes
$ws.Range("E147").Value = 391
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 14

# Mongolia (row 169): Casos activos, Recuperados
$ws.Range("D169").Value = 185
$ws.Range("E169").Value = 35
